$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-like number format on price cells whose new values look numeric,
# so Excel stores them as text (matching the original inlineStr formatting)
# instead of silently converting them to floating point numbers.
$textFormatCells = @("D4", "D5", "D6", "D7", "D8", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.844.20"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "1.733.33"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "231.19"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "0.5145"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("D8").Value = "0.2772"
$ws.Range("E8").Value = "  +4.67%  "
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").Value = "0.06105"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").Value = "1.746.64"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "0.07022"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "15.21"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "0.6416"
$ws.Range("E14").Value = "  +3.67%  "
$ws.Range("D15").Value = "4.520"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").Value = "76.78"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "25.832.42"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "11.50"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "0.000006625"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").Value = "1.966.32"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "4.143"
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("D24").Value = "8.705"
$ws.Range("E24").Value = "  +5.94%  "
$ws.Range("D25").Value = "5.116"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "139.91"
$ws.Range("E26").Value = "  +3.00%  "
$ws.Range("D27").Value = "1.515"
$ws.Range("E27").Value = "  +3.17%  "
$ws.Range("D28").Value = "15.02"
$ws.Range("D29").Value = "1.794"
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("D30").Value = "101.97"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").Value = "0.08301"
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("D32").Value = "3.691"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("E33").Value = "  +1.85%  "
$ws.Range("D34").Value = "0.04477"
$ws.Range("E34").Value = "  +2.63%  "
$ws.Range("D35").Value = "2.617"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("D36").Value = "0.9795"
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("D37").Value = "0.6131"
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("D38").Value = "2.639"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("D40").Value = "1.930"
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("D41").Value = "0.9998"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "100.41"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "0.3813"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").Value = "0.7291"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").Value = "4.958"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").Value = "0.05382"
$ws.Range("E46").Value = "  -1.84%  "
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("D48").Value = "6.234"
$ws.Range("E48").Value = "  +5.62%  "
$ws.Range("D49").Value = "52.92"
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").Value = "30.01"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").Value = "7.629"
$ws.Range("E51").Value = "  +3.46%  "
